# Criação da funcionalidade de Alteração do Destino de Alta
# Replace the "farmacia" user with "flavia" in the grant statements list
# located on the "grants por usuario" sheet (column B, rows 59-96).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grants por usuario")

for ($r = 59; $r -le 96; $r++) {
    $ws.Cells.Item($r, 2).Value = "flavia"
}

$ws.Activate()
$ws.Range("A1").Select()
